$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Program_sheet: rename a couple of program entries, add two new rows, and
# blank out the old notes cell in D14.
# ---------------------------------------------------------------------------
$progSheet = $wb.Worksheets.Item("Program_sheet")

# Row 2: rename program/description text
$progSheet.Range("A2").Value = "Yxyxyx"
$progSheet.Range("B2").Value = "ababab"

# New row 8 (previously unused row between 7 and 9)
$progSheet.Range("A8").Value = "8playwright8"
$progSheet.Range("B8").Value = "auto"
$progSheet.Range("C8").Value = "active"

# New row 12 (previously unused row between 11 and 14)
$progSheet.Range("A12").Value = "ML"
$progSheet.Range("B12").Value = "Machines"
$progSheet.Range("C12").Value = "Active"

# Row 14: clear out the old multi-line validation notes text
$progSheet.Range("D14").Value = ""

# ---------------------------------------------------------------------------
# Batch: "MobileSeleniumAuto" was renamed to "activa" everywhere it is used,
# and a few "Number of Classes" values were bumped up.
# ---------------------------------------------------------------------------
$batchSheet = $wb.Worksheets.Item("Batch")

$batchSheet.Range("B2").Value = "activa"
$batchSheet.Range("B3").Value = "activa"
$batchSheet.Range("B4").Value = "activa"
$batchSheet.Range("B5").Value = "activa"
$batchSheet.Range("B7").Value = "activa"
$batchSheet.Range("B8").Value = "activa"
$batchSheet.Range("B9").Value = "activa"
$batchSheet.Range("B10").Value = "activa"

$batchSheet.Range("C2").Value = 13
$batchSheet.Range("C4").Value = 14
$batchSheet.Range("C5").Value = 15
